$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> Dcn -> Erbb4)
$ws.Range("G2").Value = 1.006697333333333
$ws.Range("H2").Value = 3.020092
$ws.Range("I2").Value = 0.0001985651645046208
$ws.Range("J2").Value = 0.0001985651645046208
$ws.Range("Q2").Value = 0.006392863632444444
$ws.Range("R2").Value = 0.057535772692
$ws.Range("S2").Value = 0.0001985651645046208
$ws.Range("T2").Value = 0.0001985651645046208

# Row 3 (FAPs -> Dcn -> Erbb4)
$ws.Range("I3").Value = 0.9806494927176636
$ws.Range("J3").Value = 0.9806494927176637
$ws.Range("S3").Value = 0.9806494927176636
$ws.Range("T3").Value = 0.9806494927176637

# Row 4 (MuSCs -> Dcn -> Erbb4)
$ws.Range("G4").Value = 97.097641
$ws.Range("H4").Value = 291.292923
$ws.Range("I4").Value = 0.01915194211783179
$ws.Range("J4").Value = 0.01915194211783179
$ws.Range("Q4").Value = 0.6166023862303333
$ws.Range("R4").Value = 5.549421476072999
$ws.Range("S4").Value = 0.01915194211783179
$ws.Range("T4").Value = 0.01915194211783179
